# Applies the tracked-change style edits described by the commit diff:
#  1. "Team: " paragraph - collapse the "Arunsundar" proofed run back into
#     its neighbours (no text change, just a run merge / proofErr removal).
#  2. "...to the implementation; Factory and Decorator Design Patterns"
#     -> "...to the implementation; Factory and Composite Design Patterns"
#     (split into 3 runs, middle one carries the actual word change).
#  3. "Work done by Arunsundar Kannan" paragraph - same run-merge pattern
#     as (1).
#  4. "Implementing this design pattern makes the code open..." - merge
#     three runs into one (no text change).
#  5. Final sentence "...have been satisfied." - split around the
#     pre-existing _GoBack bookmark into "...have bee" | bookmark |
#     "n satisfied." (no text change).

$d = $word.ActiveDocument

# --- 1. "Team: Arunsundar Kannan, Ashish ..." -----------------------------
$d.Content.Find.Execute(
    ": Arunsundar Kannan, Ashish ", $true, $false, $false, $false, $false,
    $true, 1, $false, ": Arunsundar Kannan, Ashish ", 2)

# --- 2. Factory and Decorator -> Factory and Composite --------------------
$rng = $d.Content
$rng.Find.Execute(" to the implementation; Factory and Decorator Design Patterns")
$start = $rng.Start

$len1 = 13   # " to the imple"
$len2old = 32   # "mentation; Factory and Decorator"

$part1 = $d.Range($start, $start + $len1)
$ft1 = $part1.FormattedText
$part1.FormattedText = $ft1

$part2 = $d.Range($start + $len1, $start + $len1 + $len2old)
$ft2 = $part2.FormattedText
$ft2.Text = "mentation; Factory and Composite"
$part2.FormattedText = $ft2

# --- 3. "Work done by Arunsundar Kannan - " -------------------------------
$d.Content.Find.Execute(
    "Work done by Arunsundar Kannan " + [char]0x2013 + " ", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Work done by Arunsundar Kannan " + [char]0x2013 + " ", 2)

# --- 4. Merge the "Implementing this design pattern..." runs -------------
$d.Content.Find.Execute(
    "Implementing this design pattern makes the code open for extension and we are encapsulating what is varying ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Implementing this design pattern makes the code open for extension and we are encapsulating what is varying ",
    2)

# --- 5. Split the final sentence around the existing _GoBack bookmark ----
$rng = $d.Content
$rng.Find.Execute("The final step would be to test the project based on the use case document and verify that all the requirements have bee")
$head = $d.Range($rng.Start, $rng.End)
$splitPos = $rng.End

$tailFull = $d.Content
$tailFull.Find.Execute("n satisfied.")
$tail = $d.Range($splitPos, $tailFull.End)
$ftTail = $tail.FormattedText
$tail.FormattedText = $ftTail
